$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Remove the stray <w:bookmarkStart/bookmarkEnd w:name="_GoBack"/> that
#    trails the "...existing code." run in the first (Project Title/Purpose)
#    table. InsertXML replaces the whole enclosing paragraph in this engine,
#    so the Find anchor/replacement XML below re-states every run of that
#    paragraph verbatim and simply omits the bookmark pair at the end.
# ---------------------------------------------------------------------------
$rng = $d.Content
$wholePara1 = "Yocto+SPDX is an ongoing opensource software development project.  " + `
    "The purpose of UNOs Yocto+SPDX group is to migrate/fork and manage the current " + `
    "project to a GitHub Repository and make improvements to the existing code."
$found = $rng.Find.Execute(
    $wholePara1,
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "anchor 1 not found" }

$xml1 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document $wns><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Yocto+SPDX is an ongoing opensource software development project.  The purpose of UNOs Yocto+SPDX group is to migrate</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>/fork</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">and manage </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>the current project to a GitHub Repository and make improvements to the existing code.</w:t></w:r></w:p></w:body></w:document>
"@
$rng.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Split the paragraph that used to begin with "The Yocto+SPDX project is
#    built..." by inserting two brand-new paragraphs in front of it:
#      - a short paragraph introducing the SPDX specification
#      - a "Source: https://spdx.org/about-spdx/what-is-spdx" paragraph,
#        which also picks up the relocated _GoBack bookmark (explicitly
#        re-using id 0, freed up by step 1 above)
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "The Yocto+SPDX project is built",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "anchor 2 not found" }
$rng2.Collapse(1)   # wdCollapseStart - zero-length range right before that paragraph

$xml2 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document $wns><w:body>
<w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>The Software Package Data Exchange® (SPDX®) specification is a standard format for communicating the components, licenses and copyrights associated with a software package.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Source: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>https://spdx.org/about-spdx/what-is-spdx</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
</w:body></w:document>
"@
$rng2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) Rework the closing "Description .../Source: ..." paragraph:
#      - give its paragraph mark run-properties a sz/szCs of 20
#      - prepend a new "Yocto+SPDX " run before "Description "
#      - give every run sz/szCs of 20
#      - merge "Source: " and the bare URL into a single run's text
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute(
    "Description Source: https://spdx.org/tools/community/yoctospdx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "anchor 3 not found" }

$xml3 = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document $wns><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Yocto+SPDX </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Description </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Source: https://spdx.org/tools/community/yoctospdx</w:t></w:r></w:p></w:body></w:document>
"@
$rng3.InsertXML($xml3)
